# DaySale_2025-08-05: add two missing-stock rows (SPASMO-DIGESTIN, SUGARLO PLUS)
# that were omitted alphabetically among existing rows, add two more rows at the
# bottom (KALONA / shaving machine) that were previously missing before the
# closing total / footer bar, refresh the running total and the "printed at"
# timestamp string.

function Set-TextValue($range, [string]$text) {
    # Some columns (L, P) carry a numeric NumberFormat even though the sheet
    # stores their content as literal text (e.g. "1", "0", "25.7400").
    # Writing straight into .Value on such a cell gets silently coerced to a
    # number (losing the exact literal text / trailing zeros). Flipping the
    # format to Text for the write and then restoring it avoids that without
    # disturbing the cell's visual style.
    $origFmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = $origFmt
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert "SPASMO-DIGESTIN 30 TABS." right before "STREPTOQUIN 20 TABLETS"
#    (currently row 17), pushing it and everything below down by one row.
# ---------------------------------------------------------------------------
$ws.Range("A17:Q17").Insert()
$ws.Range("A18:Q18").Copy($ws.Range("A17:Q17"))

Set-TextValue $ws.Range("C17") "SPASMO-DIGESTIN 30 TABS."
Set-TextValue $ws.Range("H17") "2:1"
Set-TextValue $ws.Range("L17") "1"
Set-TextValue $ws.Range("N17") "78.00"
Set-TextValue $ws.Range("P17") "25.7400"
Set-TextValue $ws.Range("Q17") "0:1"

# ---------------------------------------------------------------------------
# 2) Insert "SUGARLO PLUS 50/1000MG 30 F.C. TABS" right before
#    "TRILLERG EYE DROPS 10 ML" (now row 19 after the previous insert).
# ---------------------------------------------------------------------------
$ws.Range("A19:Q19").Insert()
$ws.Range("A20:Q20").Copy($ws.Range("A19:Q19"))

Set-TextValue $ws.Range("C19") "SUGARLO PLUS 50/1000MG 30 F.C. TABS"
Set-TextValue $ws.Range("H19") "2:1"
Set-TextValue $ws.Range("L19") "1"
Set-TextValue $ws.Range("N19") "136.50"
Set-TextValue $ws.Range("P19") "45.0450"
Set-TextValue $ws.Range("Q19") "0:1"

# ---------------------------------------------------------------------------
# 3) Two more data rows were missing before the total row. They belong right
#    after "ماكينه حلاقه جليت فليكتور" (shaving machine) ... but that item IS
#    one of the two new rows; the other new row ("كالونا") is new too. Both
#    sit at the very end of the data block, just above the total/footer bars
#    (which currently sit at rows 28 and 29).
# ---------------------------------------------------------------------------
$ws.Range("A28:Q29").Insert()
$ws.Range("A27:Q27").Copy($ws.Range("A28:Q28"))
$ws.Range("A27:Q27").Copy($ws.Range("A29:Q29"))

Set-TextValue $ws.Range("C28") "كالونا "
Set-TextValue $ws.Range("H28") "0:0"
Set-TextValue $ws.Range("L28") "0"
Set-TextValue $ws.Range("N28") "15.00"
Set-TextValue $ws.Range("P28") "15.0000"
Set-TextValue $ws.Range("Q28") "1:0"

Set-TextValue $ws.Range("C29") "ماكينه حلاقه جليت فليكتور"
Set-TextValue $ws.Range("H29") "21:0"
Set-TextValue $ws.Range("L29") "0"
Set-TextValue $ws.Range("N29") "15.00"
Set-TextValue $ws.Range("P29") "15.0000"
Set-TextValue $ws.Range("Q29") "1:0"

# ---------------------------------------------------------------------------
# 4) Fix up the "م" (row index) column, which is always just (row - 6), for
#    every data row now that two rows were inserted in the middle.
# ---------------------------------------------------------------------------
for ($r = 7; $r -le 29; $r++) {
    $ws.Range("A$r").Value = $r - 6
}

# ---------------------------------------------------------------------------
# 5) Recompute the grand total of the "selling price" column (P7:P29) shown
#    on the total row (now row 30), and bump the "printed at" timestamp
#    (now row 31) by one minute.
# ---------------------------------------------------------------------------
$total = 0
for ($r = 7; $r -le 29; $r++) {
    $total += [double]$ws.Range("P$r").Text
}
$ws.Range("P30").Value = $total

Set-TextValue $ws.Range("A31") "Tuesday, 5 August, 2025 11:06 AM"
